# Add a new worksheet "Sheet1" with project feature/description data,
# and update the selection/view state on both sheets to match the
# target workbook state.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet (after the existing "Project Planner" sheet) ---
$planner0 = $wb.Worksheets.Item("Project Planner")
$newSheet = $wb.Worksheets.Add($null, $planner0)
$newSheet.Name = "Sheet1"

# --- Populate the new sheet with feature/description data ---
# (values are entered in the same order the original author typed them,
# so shared-string allocation order matches the target file)
$newSheet.Range("A1").Value = "Features"
$newSheet.Range("C1").Value = "Schedule Completion Week"
$newSheet.Range("D1").Value = "Actual Completion Week"

$newSheet.Range("A2").Value = "A Map base website"

$newSheet.Range("B1").Value = "Description"
$newSheet.Range("B2").Value = "A map engine is included in a website. Allow user to simply move around the map"

$newSheet.Range("A3").Value = "Basic button include ""Home"", ""News"",  ""Contact"", ""About"""
$newSheet.Range("B3").Value = "Display general information"

$newSheet.Range("A4").Value = "Allow user to search on the map base on their input"
$newSheet.Range("B4").Value = "A search box on the map allow the user to input the map area they would like to see"

$newSheet.Range("A5").Value = "Integrate above function on one website"

$newSheet.Range("A6").Value = "Side bar available"
$newSheet.Range("B6").Value = "Including some dummy button to have a feel of the user interface"

$newSheet.Range("A7").Value = "Switch between differernt base map"
$newSheet.Range("B7").Value = "Switch between solar/wind/water"

$newSheet.Range("A8").Value = "Switch between differernt map data"
$newSheet.Range("B8").Value = "Different map layer visualized for differernt energy. Example: wind at 10m, wind at 50m."

# --- Update view/selection state ---
# "Project Planner" sheet: selection moved, no longer the active tab.
$planner = $wb.Worksheets.Item("Project Planner")
$planner.Range("B22").Select() | Out-Null

# Make the new sheet the active (selected) tab, matching tabSelected="1" on
# the new sheet in the target workbook.
$newSheet.Activate()
$newSheet.Range("A9").Select() | Out-Null
